$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data for 12/12/2025 (serial 46003)
$ws.Range("A46").Value = 46003
$ws.Range("B46").Value = 543
$ws.Range("C46").Value = 18
$ws.Range("D46").Value = 525

# Update the selection to reflect the newly active row, as seen in the saved file
$ws.Range("A46:D46").Select()
